$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete columns that are no longer needed (shifts remaining columns left).
# Delete from right to left so earlier column letters stay valid.
# Columns to delete (by their ORIGINAL letters): C, E, G, I, T, U, V, W, X, Y, Z, AA
$ws.Range("T1:AA1").EntireColumn.Delete() | Out-Null
$ws.Range("I1").EntireColumn.Delete() | Out-Null
$ws.Range("G1").EntireColumn.Delete() | Out-Null
$ws.Range("E1").EntireColumn.Delete() | Out-Null
$ws.Range("C1").EntireColumn.Delete() | Out-Null

# After the deletions, the remaining headers (A..O) are:
# A l.p | B Data dodania | C Link | D Adres | E podzielnica | F Posrednik? |
# G Telefon | H Cena | I Cena/m2 | J m2 | K Pokoje | L Pietro |
# M Rodzaj mieszkania | N Notatka | O Tabela dzwonienie
#
# Target order needs Notatka and "Tabela dzwonienie" swapped, with a new
# "Max Cena Kupna" column inserted between them. Insert a blank column at N
# first (this shifts Notatka -> O and Tabela dzwonienie -> P), then fill in
# the three header cells with their final text.
$ws.Range("N1").EntireColumn.Insert() | Out-Null

$ws.Range("N1").Value = "Tabela dzwonienie"
$ws.Range("O1").Value = "Max Cena Kupna"
$ws.Range("P1").Value = "Notatka"
